# Nieuwe data toegevoegd via Streamlit op 2024-12-03 18:09:10
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 71

$ws.Cells.Item($row, 1).Value = "CompaNanny"
$ws.Cells.Item($row, 2).Value = "CompaNanny Statenkwartier BSO"
$ws.Cells.Item($row, 3).Value = "VGO"

# The report date column holds plain text dates (e.g. "2023-03-28") rather
# than real date serials elsewhere in this sheet, so force text formatting
# before assigning the value to stop it being auto-parsed as a date, then
# drop the formatting again so the cell keeps the workbook's default style.
$dateCell = $ws.Cells.Item($row, 4)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2023-03-28"
$dateCell.ClearFormats()

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
